{"js": "// Sequential 1:1 text replacements for every paragraph in the document body,\n// in document order: paragraph 0 is the date line above the table, and\n// paragraphs 1..100 are the 100 table-cell paragraphs (20 rows x 5 cols).\nconst pairs = [\n  [\"2025-12-11 Thursday\", \"2025-12-12 Friday\"],\n  [\"17+4=\", \"43+29=\"],\n  [\"77-31=\", \"55-26=\"],\n  [\"59+0=\", \"64+7=\"],\n  [\"82-12=\", \"83-69=\"],\n  [\"98-93=\", \"87-35=\"],\n  [\"65-16=\", \"54+31=\"],\n  [\"60-32=\", \"98-85=\"],\n  [\"26+53=\", \"97-18=\"],\n  [\"75-28=\", \"83+2=\"],\n  [\"46+39=\", \"83-59=\"],\n  [\"44-2=\", \"56+21=\"],\n  [\"70-65=\", \"81+13=\"],\n  [\"43+47=\", \"68-12=\"],\n  [\"6+76=\", \"77+22=\"],\n  [\"51-43=\", \"60-59=\"],\n  [\"42+31=\", \"53-5=\"],\n  [\"96-92=\", \"75-64=\"],\n  [\"78-30=\", \"71-41=\"],\n  [\"17+70=\", \"42+22=\"],\n  [\"1+61=\", \"46+1=\"],\n  [\"5+17=\", \"29+63=\"],\n  [\"80-4=\", \"99-5=\"],\n  [\"10+68=\", \"28+63=\"],\n  [\"75-25=\", \"40+27=\"],\n  [\"99-50=\", \"96-46=\"],\n  [\"37+28=\", \"85-54=\"],\n  [\"82-81=\", \"34-1=\"],\n  [\"93-51=\", \"20+78=\"],\n  [\"52+32=\", \"3+0=\"],\n  [\"75+20=\", \"7+70=\"],\n  [\"51+38=\", \"47-10=\"],\n  [\"34+30=\", \"97-58=\"],\n  [\"97-39=\", \"68+6=\"],\n  [\"27+38=\", \"57-4=\"],\n  [\"57-0=\", \"34+33=\"],\n  [\"1+30=\", \"65+3=\"],\n  [\"67-7=\", \"82+16=\"],\n  [\"12+24=\", \"82-50=\"],\n  [\"9+83=\", \"11+45=\"],\n  [\"30+0=\", \"52-17=\"],\n  [\"96-63=\", \"14+33=\"],\n  [\"66-65=\", \"13+41=\"],\n  [\"93-54=\", \"74-41=\"],\n  [\"91-15=\", \"36+30=\"],\n  [\"96-21=\", \"24+68=\"],\n  [\"3+11=\", \"10+34=\"],\n  [\"76+11=\", \"29+67=\"],\n  [\"26-11=\", \"54-48=\"],\n  [\"64+13=\", \"20+1=\"],\n  [\"41+5=\", \"32+57=\"],\n  [\"42+31=\", \"41+9=\"],\n  [\"74-27=\", \"45+54=\"],\n  [\"29+0=\", \"60+14=\"],\n  [\"66-48=\", \"84-58=\"],\n  [\"11+10=\", \"75+17=\"],\n  [\"79-64=\", \"71-63=\"],\n  [\"34+54=\", \"86-12=\"],\n  [\"40-0=\", \"25+19=\"],\n  [\"57-11=\", \"31+28=\"],\n  [\"14-11=\", \"6+89=\"],\n  [\"31+2=\", \"97-19=\"],\n  [\"70-64=\", \"13+20=\"],\n  [\"61-21=\", \"70-37=\"],\n  [\"8+19=\", \"9+15=\"],\n  [\"46-35=\", \"60-40=\"],\n  [\"80-73=\", \"67-26=\"],\n  [\"5+69=\", \"67-49=\"],\n  [\"22+1=\", \"29-17=\"],\n  [\"94-64=\", \"70+7=\"],\n  [\"60+35=\", \"64-32=\"],\n  [\"9+19=\", \"72-29=\"],\n  [\"82-65=\", \"38+6=\"],\n  [\"95-83=\", \"85+3=\"],\n  [\"66-21=\", \"62-53=\"],\n  [\"11+72=\", \"18+54=\"],\n  [\"17+73=\", \"49+47=\"],\n  [\"12+56=\", \"25+53=\"],\n  [\"76+12=\", \"43-40=\"],\n  [\"40+12=\", \"77-19=\"],\n  [\"56-54=\", \"78-12=\"],\n  [\"26+8=\", \"45-39=\"],\n  [\"47+6=\", \"99-2=\"],\n  [\"58+2=\", \"13+64=\"],\n  [\"19-13=\", \"12+83=\"],\n  [\"26+49=\", \"59-12=\"],\n  [\"89-62=\", \"25-3=\"],\n  [\"2+61=\", \"58+11=\"],\n  [\"92-64=\", \"80-54=\"],\n  [\"93-56=\", \"46+52=\"],\n  [\"20+36=\", \"95-75=\"],\n  [\"20+13=\", \"87-3=\"],\n  [\"46-40=\", \"0+3=\"],\n  [\"27-18=\", \"81-46=\"],\n  [\"87-40=\", \"24+20=\"],\n  [\"56-44=\", \"22-5=\"],\n  [\"21-14=\", \"55+36=\"],\n  [\"64-17=\", \"82-6=\"],\n  [\"38-9=\", \"82+4=\"],\n  [\"27+70=\", \"39-32=\"],\n  [\"96-56=\", \"51-24=\"]\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nif (paragraphs.items.length !== pairs.length) {\n  throw new Error(\n    \"Paragraph count (\" + paragraphs.items.length +\n    \") does not match expected pair count (\" + pairs.length + \")\"\n  );\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const [oldText, newText] = pairs[i];\n  if (para.text !== oldText) {\n    throw new Error(\n      \"Paragraph \" + i + \" text mismatch: expected \\\"\" + oldText +\n      \"\\\" but found \\\"\" + para.text + \"\\\"\"\n    );\n  }\n  if (oldText !== newText) {\n    para.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$dateOld = '2025-12-11 Thursday'\n$dateNew = '2025-12-12 Friday'\n\n# Row-major list of [old, new] pairs for the 20x5 table (100 cells)\n$cellPairs = @(\n    @('17+4=', '43+29='),\n    @('77-31=', '55-26='),\n    @('59+0=', '64+7='),\n    @('82-12=', '83-69='),\n    @('98-93=', '87-35='),\n    @('65-16=', '54+31='),\n    @('60-32=', '98-85='),\n    @('26+53=', '97-18='),\n    @('75-28=', '83+2='),\n    @('46+39=', '83-59='),\n    @('44-2=', '56+21='),\n    @('70-65=', '81+13='),\n    @('43+47=', '68-12='),\n    @('6+76=', '77+22='),\n    @('51-43=', '60-59='),\n    @('42+31=', '53-5='),\n    @('96-92=', '75-64='),\n    @('78-30=', '71-41='),\n    @('17+70=', '42+22='),\n    @('1+61=', '46+1='),\n    @('5+17=', '29+63='),\n    @('80-4=', '99-5='),\n    @('10+68=', '28+63='),\n    @('75-25=', '40+27='),\n    @('99-50=', '96-46='),\n    @('37+28=', '85-54='),\n    @('82-81=', '34-1='),\n    @('93-51=', '20+78='),\n    @('52+32=', '3+0='),\n    @('75+20=', '7+70='),\n    @('51+38=', '47-10='),\n    @('34+30=', '97-58='),\n    @('97-39=', '68+6='),\n    @('27+38=', '57-4='),\n    @('57-0=', '34+33='),\n    @('1+30=', '65+3='),\n    @('67-7=', '82+16='),\n    @('12+24=', '82-50='),\n    @('9+83=', '11+45='),\n    @('30+0=', '52-17='),\n    @('96-63=', '14+33='),\n    @('66-65=', '13+41='),\n    @('93-54=', '74-41='),\n    @('91-15=', '36+30='),\n    @('96-21=', '24+68='),\n    @('3+11=', '10+34='),\n    @('76+11=', '29+67='),\n    @('26-11=', '54-48='),\n    @('64+13=', '20+1='),\n    @('41+5=', '32+57='),\n    @('42+31=', '41+9='),\n    @('74-27=', '45+54='),\n    @('29+0=', '60+14='),\n    @('66-48=', '84-58='),\n    @('11+10=', '75+17='),\n    @('79-64=', '71-63='),\n    @('34+54=', '86-12='),\n    @('40-0=', '25+19='),\n    @('57-11=', '31+28='),\n    @('14-11=', '6+89='),\n    @('31+2=', '97-19='),\n    @('70-64=', '13+20='),\n    @('61-21=', '70-37='),\n    @('8+19=', '9+15='),\n    @('46-35=', '60-40='),\n    @('80-73=', '67-26='),\n    @('5+69=', '67-49='),\n    @('22+1=', '29-17='),\n    @('94-64=', '70+7='),\n    @('60+35=', '64-32='),\n    @('9+19=', '72-29='),\n    @('82-65=', '38+6='),\n    @('95-83=', '85+3='),\n    @('66-21=', '62-53='),\n    @('11+72=', '18+54='),\n    @('17+73=', '49+47='),\n    @('12+56=', '25+53='),\n    @('76+12=', '43-40='),\n    @('40+12=', '77-19='),\n    @('56-54=', '78-12='),\n    @('26+8=', '45-39='),\n    @('47+6=', '99-2='),\n    @('58+2=', '13+64='),\n    @('19-13=', '12+83='),\n    @('26+49=', '59-12='),\n    @('89-62=', '25-3='),\n    @('2+61=', '58+11='),\n    @('92-64=', '80-54='),\n    @('93-56=', '46+52='),\n    @('20+36=', '95-75='),\n    @('20+13=', '87-3='),\n    @('46-40=', '0+3='),\n    @('27-18=', '81-46='),\n    @('87-40=', '24+20='),\n    @('56-44=', '22-5='),\n    @('21-14=', '55+36='),\n    @('64-17=', '82-6='),\n    @('38-9=', '82+4='),\n    @('27+70=', '39-32='),\n    @('96-56=', '51-24=')\n)\n\n# --- Update the date line (first paragraph, above the table) ---\n$titlePara = $d.Paragraphs.Item(1)\n$titleRange = $titlePara.Range\n# Trim the trailing paragraph mark from the comparison text\n$titleText = $titleRange.Text.TrimEnd([char]13, [char]7)\nif ($titleText -ne $dateOld) {\n    throw \"Title paragraph text mismatch: expected '$dateOld' but found '$titleText'\"\n}\n$titleRange.Text = $dateNew\n\n# --- Update every cell of the first (and only) table, row-major, 20 rows x 5 cols ---\n$tbl = $d.Tables.Item(1)\n$rows = $tbl.Rows.Count\n$cols = $tbl.Columns.Count\nif ($rows * $cols -ne $cellPairs.Count) {\n    throw \"Table size ($rows x $cols = $($rows * $cols)) does not match expected cell count ($($cellPairs.Count))\"\n}\n\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $pair = $cellPairs[$idx]\n        $old = $pair[0]\n        $new = $pair[1]\n\n        $cell = $tbl.Cell($r, $c)\n        $cellRange = $cell.Range\n        $cellText = $cellRange.Text.TrimEnd([char]13, [char]7)\n        if ($cellText -ne $old) {\n            throw \"Cell ($r,$c) text mismatch: expected '$old' but found '$cellText'\"\n        }\n        if ($old -ne $new) {\n            $cellRange.Text = $new\n        }\n        $idx++\n    }\n}\n\nWrite-Output \"Updated $idx cells plus the title line.\"\n"}
